$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "EOF" marker from A33 down to A34, and put the "EOC" marker into B33
# (engine position block gets its own EOC separator, and wing position section's EOF moves down)
$ws.Range("A33").ClearContents()
$ws.Range("B33").Value = "EOC"
$ws.Range("A34").Value = "EOF"

# Update the active selection/view as recorded in the saved workbook
$ws.Range("B33").Select()
$excel.ActiveWindow.ScrollRow = 13
